$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "Compass Error Compass Error Compass disconnected ."
$ws.Range("C2").Value = "Compass Error"
$ws.Range("D2").Value = "0-1"

# Row 3
$ws.Range("A3").Value = 0
$ws.Range("C3").Value = "Compass disconnected"
$ws.Range("D3").Value = "4-5"

# Row 4
$ws.Range("A4").Value = 0
$ws.Range("C4").Value = "Compass Error"
$ws.Range("D4").Value = "2-3"

# Row 5
$ws.Range("B5").Value = "Critical low battery Aircraft in Auto Power Off Protection Forced landing in progress ."
$ws.Range("C5").Value = "Critical low battery"
$ws.Range("D5").Value = "0-2"
